$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13; existing rows 13-17 shift down to 14-18.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 45233
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = 100112030
$ws.Range("G13").Value = "Poroto granado"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 1050
$ws.Range("K13").Value = 1400
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1438
$ws.Range("N13").Value = "`$/kilo"
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 1438
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"

$wb.Save()
